$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registros")

# Row 5 (A5:C5 "test" entry) had trailing empty placeholder cells D5:K5 -
# drop them so the row only keeps its real data (A5:C5).
$ws.Range("D5:K5").ClearContents()

# New collaborator record appended as row 6. Force text format first so
# Excel keeps the values as plain strings (matching the rest of the
# sheet) instead of auto-coercing the date-looking / numeric-looking
# strings into real dates or numbers.
$newRow = $ws.Range("A6:K6")
$newRow.NumberFormat = "@"

$ws.Range("A6").Value = "01/01/2023"
$ws.Range("B6").Value = "01/01/2023"
$ws.Range("C6").Value = "LUIS SANTANA"
$ws.Range("D6").Value = "99"
$ws.Range("E6").Value = "9"
$ws.Range("F6").Value = "9"
$ws.Range("G6").Value = "9"
$ws.Range("H6").Value = "90"
$ws.Range("I6").Value = "00"
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
